# Apply cryptos list price/volume refresh as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.904.40"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.807.48"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.32"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4448"
$ws.Range("E7").Value = "  +5.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3676"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07341"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8562"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "1.903.01"
$ws.Range("E12").Value = "  +4.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.619"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.51"
$ws.Range("E14").Value = "  +3.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.307"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07056"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008731"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.86"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "26.923.58"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.152"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.995"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.91"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.52"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.179"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.211"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.59"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08828"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7497"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.174"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.938"
$ws.Range("E33").Value = "  +4.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.458"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9997"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.087"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01966"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05189"
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5319"
$ws.Range("E39").Value = "  +5.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.865"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.013"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1691"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5157"
$ws.Range("E43").Value = "  +8.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.417"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.989"
$ws.Range("E45").Value = "  +6.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.55"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.35"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.667"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06323"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9170"
$ws.Range("E51").Value = "  +0.45%  "
